# Clean up code and fix output
# Adds a new "Yearly demand" worksheet at the end of the workbook with
# hourly net-demand data for three representative days.

$wb = $excel.ActiveWorkbook
$originalActiveSheet = $wb.ActiveSheet

# Duplicate the last sheet (placing the copy right after it, i.e. at the
# very end of the workbook) so the new sheet inherits the same sheet
# defaults (page setup, outline options, etc.) as the rest of the model,
# then wipe it completely before filling it with the new content.
$templateSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$templateSheet.Copy([System.Reflection.Missing]::Value, $templateSheet)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Yearly demand"
$ws.Cells.Clear()

# Header row: hours 0-23 across columns B:Y
for ($col = 2; $col -le 25; $col++) {
    $ws.Cells.Item(1, $col).Value = $col - 2
}

# Day index column A, rows 2-4 -> 0,1,2
for ($row = 2; $row -le 4; $row++) {
    $ws.Cells.Item($row, 1).Value = $row - 2
}

$row2 = @(-32.5,-19.5,-13,-13,-13,142.5,291.5,327,388.5,502,596,670.5,745,651,576.5,502,320.5,139,32,-117,-97.5,-78,-52,-39)
$row3 = @(-32.5,-19.5,-13,0,0,-19.5,0,324,486,648,729,751.5,583,567,333.5,340,243,57.99999999999999,-130,0,0,-78,0,-39)
$row4 = @(-32.5,-19.5,0,0,0,-19.5,0,0,81,324,567,589.5,648,567,324,162,81,0,-130,0,0,0,0,-39)

$data = @($row2, $row3, $row4)

for ($i = 0; $i -lt 3; $i++) {
    $r = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt 24; $j++) {
        $ws.Cells.Item($r, $j + 2).Value = $vals[$j]
    }
}

# Styling to match the rest of the workbook (e.g. "Connected Households"):
# header row (B1:Y1) and the day-index column (A2:A4) use the bold,
# centered/top-aligned, thin-bordered style already used elsewhere in the
# workbook. Copy that format instead of building it from scratch so the
# existing style entry is reused rather than creating a near-duplicate.
$styleSource = $wb.Worksheets.Item("Connected Households").Range("B1")
$styleSource.Copy()
$ws.Range("B1:Y1").PasteSpecial(-4122)
$ws.Range("A2:A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A1").Select() | Out-Null

# Restore the workbook's original active sheet/selection so this purely
# additive change doesn't also shift which tab is shown on open.
$originalActiveSheet.Activate() | Out-Null
